$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 8 (pushes existing rows 8..28 down to 9..29),
# copying row 7's formatting/merges so the new row matches the table style.
$ws.Rows.Item(8).Insert()
$ws.Range("A7:Q7").Copy($ws.Range("A8:Q8"))

# Fill in the new item's data (BLOKATENS), matching the other data rows' text-based values.
$ws.Range("A8").Value = 2
$ws.Range("C8").Value = "BLOKATENS 5/80MG 28 F.C. TAB"
$ws.Range("H8").Value = "0:1"
$ws.Range("L8").Value = "1"
$ws.Range("N8").Value = "122.00"
$ws.Range("P8").Value = "61.0000"
$ws.Range("Q8").Value = "0:1"

# Update the grand-total (was row 27, now row 28 after the insert) to include the new row.
$ws.Range("P28").Value = 501.43

# Update the generated timestamp in the footer (was row 28, now row 29 after the insert).
$ws.Range("A29").Value = "Sunday, 10 August, 2025 11:29 AM"
